# Apply "fixed no fault data" edits to Sheet1 of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated data values for rows 2-11: Column B -> 100, and new E/F/G values.
$data = @(
    @{ Row = 2;  B = 100; E = 942;  F = 19703384;  G = 92519041 }
    @{ Row = 3;  B = 100; E = 1551; F = 35722808;  G = 94016619 }
    @{ Row = 4;  B = 100; E = 2154; F = 56203784;  G = 98854219 }
    @{ Row = 5;  B = 100; E = 2633; F = 73843544;  G = 103817289 }
    @{ Row = 6;  B = 100; E = 3243; F = 96230776;  G = 111838396 }
    @{ Row = 7;  B = 100; E = 3706; F = 113592304; G = 133722718 }
    @{ Row = 8;  B = 100; E = 4227; F = 128964976; G = 188760946 }
    @{ Row = 9;  B = 100; E = 4788; F = 149534432; G = 170682100 }
    @{ Row = 10; B = 100; E = 5174; F = 164975200; G = 184293197 }
    @{ Row = 11; B = 100; E = 5894; F = 186626800; G = 277777750 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B   # Column B
    $ws.Cells.Item($r, 5).Value = $entry.E   # Column E
    $ws.Cells.Item($r, 6).Value = $entry.F   # Column F
    $ws.Cells.Item($r, 7).Value = $entry.G   # Column G
}

# Widen column F to match column G's custom (bestFit, ~12.5 chars) width,
# mirroring the <col min="6" max="7".../> change in the saved XML. Column G
# already carries the bestFit width and is left untouched.
$ws.Range("F1").EntireColumn.ColumnWidth = 11.86

# Update the active selection on the sheet to I11 (single cell).
$ws.Range("I11").Select()
